# Scheduled market-data refresh: update currentAveragePrice / Leve profit
# columns (H:N) for the affected Leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 448.85715
$ws.Range("I5").Value = 673.2222
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 673.2222
$ws.Range("L5").Value = 45
$ws.Range("M5").Value = -558.2222
$ws.Range("N5").Value = -275
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 9688.200000000001
$ws.Range("I64").Value = 8146
$ws.Range("K64").Value = 8146
$ws.Range("M64").Value = -7898
# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 9688.200000000001
$ws.Range("I67").Value = 8146
$ws.Range("K67").Value = 8146
$ws.Range("M67").Value = -7288
# Row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 7412
$ws.Range("I74").Value = 5418.231
$ws.Range("K74").Value = 5418.231
$ws.Range("M74").Value = -4482.231
# Row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 7412
$ws.Range("I77").Value = 5418.231
$ws.Range("K77").Value = 27091.155
$ws.Range("M77").Value = -22411.155
# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 1377.75
$ws.Range("I96").Value = 1647.9166
$ws.Range("J96").Value = 567.25
$ws.Range("K96").Value = 4943.7498
$ws.Range("L96").Value = 1701.75
$ws.Range("M96").Value = -3570.7498
$ws.Range("N96").Value = -4447.75
# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 6167.4
$ws.Range("I100").Value = 5116.1665
$ws.Range("K100").Value = 5116.1665
$ws.Range("M100").Value = -4575.1665
# Row 105 (Leve Item ID 18668)
$ws.Range("H105").Value = 14000
$ws.Range("J105").Value = 14000
$ws.Range("L105").Value = 14000
$ws.Range("N105").Value = -20988
# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 11818.4
$ws.Range("I106").Value = 4894.75
$ws.Range("J106").Value = 16434.166
$ws.Range("K106").Value = 4894.75
$ws.Range("L106").Value = 16434.166
$ws.Range("M106").Value = -4263.75
$ws.Range("N106").Value = -17696.166
# Row 118 (Leve Item ID 27958)
$ws.Range("H118").Value = 1388.8572
$ws.Range("I118").Value = 313.5
$ws.Range("K118").Value = 940.5
$ws.Range("M118").Value = 716.5
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1535.8928
$ws.Range("I132").Value = 1569.4231
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 4708.2693
$ws.Range("L132").Value = 3300
$ws.Range("M132").Value = -2178.2693
$ws.Range("N132").Value = -8360

$ws = $wb.Worksheets.Item("ARM")
# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 273.77777
$ws.Range("I5").Value = 219.66667
$ws.Range("J5").Value = 382
$ws.Range("K5").Value = 219.66667
$ws.Range("L5").Value = 382
$ws.Range("M5").Value = -107.66667
$ws.Range("N5").Value = -606
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 1403.6974
$ws.Range("I32").Value = 1058.452
$ws.Range("K32").Value = 1058.452
$ws.Range("M32").Value = -771.452
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4698.108
$ws.Range("I61").Value = 4759.7354
$ws.Range("J61").Value = 3999.6667
$ws.Range("K61").Value = 4759.7354
$ws.Range("L61").Value = 3999.6667
$ws.Range("M61").Value = -4547.7354
$ws.Range("N61").Value = -4423.6667
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 1241.5264
$ws.Range("I97").Value = 1603.6923
$ws.Range("K97").Value = 1603.6923
$ws.Range("M97").Value = -1107.6923
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2288.2222
$ws.Range("I102").Value = 2085
$ws.Range("K102").Value = 2085
$ws.Range("M102").Value = -463
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3701.1
$ws.Range("I122").Value = 1799.4
$ws.Range("K122").Value = 5398.200000000001
$ws.Range("M122").Value = -2948.200000000001
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2157.0156
$ws.Range("I132").Value = 1483.6
$ws.Range("J132").Value = 6272.3335
$ws.Range("K132").Value = 4450.799999999999
$ws.Range("L132").Value = 18817.0005
$ws.Range("M132").Value = -1920.799999999999
$ws.Range("N132").Value = -23877.0005
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4698.108
$ws.Range("I136").Value = 4759.7354
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 14279.2062
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = -11729.2062
$ws.Range("N136").Value = -17099.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 273.77777
$ws.Range("I4").Value = 219.66667
$ws.Range("J4").Value = 382
$ws.Range("K4").Value = 219.66667
$ws.Range("L4").Value = 382
$ws.Range("M4").Value = -104.66667
$ws.Range("N4").Value = -612
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 2362.5557
$ws.Range("I22").Value = 741.25
$ws.Range("J22").Value = 3659.6
$ws.Range("K22").Value = 741.25
$ws.Range("L22").Value = 3659.6
$ws.Range("M22").Value = -568.25
$ws.Range("N22").Value = -4005.6
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 564.35297
$ws.Range("J94").Value = 572.3333
$ws.Range("L94").Value = 572.3333
$ws.Range("N94").Value = -1474.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 46893.383
$ws.Range("I31").Value = 6765.3335
$ws.Range("K31").Value = 6765.3335
$ws.Range("M31").Value = -6470.3335
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 46893.383
$ws.Range("I34").Value = 6765.3335
$ws.Range("K34").Value = 6765.3335
$ws.Range("M34").Value = -6563.3335
# Row 82 (Leve Item ID 10799)
$ws.Range("H82").Value = 38450
# Row 85 (Leve Item ID 10799)
$ws.Range("H85").Value = 38450
# Row 88 (Leve Item ID 10608)
$ws.Range("H88").Value = 13750
$ws.Range("J88").Value = 13750
$ws.Range("L88").Value = 13750
$ws.Range("N88").Value = -14562
# Row 91 (Leve Item ID 10608)
$ws.Range("H91").Value = 13750
$ws.Range("J91").Value = 13750
$ws.Range("L91").Value = 13750
$ws.Range("N91").Value = -16558
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2374.2104
$ws.Range("I132").Value = 1449.7778
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 4349.3334
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -1819.3334
$ws.Range("N132").Value = -62102
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3825
$ws.Range("I134").Value = 2217.0908
$ws.Range("J134").Value = 7362.4
$ws.Range("K134").Value = 6651.2724
$ws.Range("L134").Value = 22087.2
$ws.Range("M134").Value = -4116.2724
$ws.Range("N134").Value = -27157.2

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 2541.75
$ws.Range("J34").Value = 5498
$ws.Range("L34").Value = 16494
$ws.Range("N34").Value = -16662
# Row 37 (Leve Item ID 9516)
$ws.Range("H37").Value = 160714.28
$ws.Range("J37").Value = 160714.28
$ws.Range("L37").Value = 482142.84
$ws.Range("N37").Value = -482366.84
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1173
$ws.Range("J68").Value = 664
$ws.Range("L68").Value = 1992
$ws.Range("N68").Value = -3614
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1173
$ws.Range("J71").Value = 664
$ws.Range("L71").Value = 5976
$ws.Range("N71").Value = -14088
# Row 116 (Leve Item ID 27866)
$ws.Range("H116").Value = 5483
$ws.Range("I116").Value = 794.3333
$ws.Range("K116").Value = 2382.9999
$ws.Range("M116").Value = 1059.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 9068.727999999999
$ws.Range("J80").Value = 11200.4
$ws.Range("L80").Value = 11200.4
$ws.Range("N80").Value = -13196.4
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 9068.727999999999
$ws.Range("J83").Value = 11200.4
$ws.Range("L83").Value = 56002
$ws.Range("N83").Value = -65986
# Row 104 (Leve Item ID 18666)
$ws.Range("H104").Value = 27600
$ws.Range("J104").Value = 27600
$ws.Range("L104").Value = 27600
$ws.Range("N104").Value = -34588
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3422.2778
$ws.Range("I126").Value = 2518.25
$ws.Range("J126").Value = 10654.5
$ws.Range("K126").Value = 7554.75
$ws.Range("L126").Value = 31963.5
$ws.Range("M126").Value = -5084.75
$ws.Range("N126").Value = -36903.5
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 4316.4443
$ws.Range("I132").Value = 3545.8125
$ws.Range("J132").Value = 10481.5
$ws.Range("K132").Value = 10637.4375
$ws.Range("L132").Value = 31444.5
$ws.Range("M132").Value = -8107.4375
$ws.Range("N132").Value = -36504.5
# Row 133 (Leve Item ID 41854)
$ws.Range("H133").Value = 69995.42999999999
$ws.Range("J133").Value = 69995.42999999999
$ws.Range("L133").Value = 69995.42999999999
$ws.Range("N133").Value = -80115.42999999999
# Row 135 (Leve Item ID 42006)
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
# Row 137 (Leve Item ID 43226)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 4856.148
$ws.Range("I22").Value = 1897.5555
$ws.Range("J22").Value = 6335.4443
$ws.Range("K22").Value = 1897.5555
$ws.Range("L22").Value = 6335.4443
$ws.Range("M22").Value = -1602.5555
$ws.Range("N22").Value = -6925.4443
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 4856.148
$ws.Range("I27").Value = 1897.5555
$ws.Range("J27").Value = 6335.4443
$ws.Range("K27").Value = 1897.5555
$ws.Range("L27").Value = 6335.4443
$ws.Range("M27").Value = -1790.5555
$ws.Range("N27").Value = -6549.4443
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 1668134.6
$ws.Range("J55").Value = 3133.25
$ws.Range("L55").Value = 3133.25
$ws.Range("N55").Value = -3479.25
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 5372.516
$ws.Range("I136").Value = 1406.9445
$ws.Range("J136").Value = 10863.308
$ws.Range("K136").Value = 4220.833500000001
$ws.Range("L136").Value = 32589.924
$ws.Range("M136").Value = -1670.833500000001
$ws.Range("N136").Value = -37689.924

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 5479.983
$ws.Range("I132").Value = 3642.0417
$ws.Range("K132").Value = 10926.1251
$ws.Range("M132").Value = -8396.125100000001
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2325.1924
$ws.Range("I136").Value = 1658
$ws.Range("K136").Value = 4974
$ws.Range("M136").Value = -2424
